$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.15
$ws.Range("B6").Value = 6.140999999999999
$ws.Range("B7").Value = 6.031999999999999
$ws.Range("D7").Value = -7.835000000000001
$ws.Range("B8").Value = 5.960000000000001
$ws.Range("D11").Value = -7.348000000000001
$ws.Range("D12").Value = -7.181
$ws.Range("D15").Value = -8.271000000000001
$ws.Range("B16").Value = 5.352
$ws.Range("B20").Value = 8.204000000000001
$ws.Range("D20").Value = -7.957000000000001
$ws.Range("B21").Value = 9.236000000000001
$ws.Range("D21").Value = -7.98
$ws.Range("D22").Value = -7.720000000000001
$ws.Range("D23").Value = -7.869999999999999
$ws.Range("B28").Value = 6.170999999999999
$ws.Range("B29").Value = 5.215999999999999
$ws.Range("D29").Value = -6.931
$ws.Range("B30").Value = 6.117
$ws.Range("B32").Value = 7.157999999999999
$ws.Range("D34").Value = -8.067
$ws.Range("B40").Value = 9.178999999999998
$ws.Range("D42").Value = -8.107000000000001
$ws.Range("D43").Value = -8.084
$ws.Range("D44").Value = -8.081
$ws.Range("D45").Value = -7.525000000000001
$ws.Range("B46").Value = 5.971
$ws.Range("D46").Value = -8.372000000000002
$ws.Range("D50").Value = -7.952
$ws.Range("B51").Value = 4.862
$ws.Range("D51").Value = -8.103
$ws.Range("B52").Value = 6.116000000000001
$ws.Range("B57").Value = 5.023
$ws.Range("D57").Value = -7.924000000000001
$ws.Range("B59").Value = 4.824999999999999
$ws.Range("B62").Value = 5.261
$ws.Range("D65").Value = -7.784999999999999
$ws.Range("B66").Value = 5.819
$ws.Range("D66").Value = -7.347
$ws.Range("D67").Value = -7.203999999999999
$ws.Range("B73").Value = 6.703
$ws.Range("B74").Value = 8.943999999999999
$ws.Range("B77").Value = 5.750999999999999
$ws.Range("D79").Value = -7.755
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("D87").Value = -8.129
$ws.Range("B92").Value = 5.371
$ws.Range("D92").Value = -6.572999999999999
$ws.Range("D97").Value = -8.151999999999999
$ws.Range("B100").Value = 6.031000000000001
